$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Divide values in A3:B25 by 10 (rows 3 through 25, columns A and B)
for ($r = 3; $r -le 25; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $bCell = $ws.Cells.Item($r, 2)
    $aVal = $aCell.Value()
    $bVal = $bCell.Value()
    $aCell.Value = $aVal / 10
    $bCell.Value = $bVal / 10
}

# Update the selected cell from F10 to F8
$ws.Range("F8").Select()
